# Apply scheduled data refresh to price/profit columns (H:N) across all class sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(12, 8).Value = 265
$ws.Cells.Item(12, 9).Value = 265
$ws.Cells.Item(12, 11).Value = 265
$ws.Cells.Item(12, 13).Value = -95
$ws.Cells.Item(40, 8).Value = 2750
$ws.Cells.Item(40, 9).Value = 1750
$ws.Cells.Item(40, 11).Value = 1750
$ws.Cells.Item(40, 13).Value = -1575
$ws.Cells.Item(51, 8).Value = 64137.945
$ws.Cells.Item(51, 10).Value = 83907.234
$ws.Cells.Item(51, 12).Value = 83907.234
$ws.Cells.Item(51, 14).Value = -84875.234
$ws.Cells.Item(112, 8).Value = 1691.0834
$ws.Cells.Item(112, 10).Value = 1915.8889
$ws.Cells.Item(112, 12).Value = 5747.6667
$ws.Cells.Item(112, 14).Value = -7963.6667
$ws.Cells.Item(137, 8).Value = 2514.2
$ws.Cells.Item(137, 9).Value = 2392.75
$ws.Cells.Item(137, 11).Value = 7178.25
$ws.Cells.Item(137, 13).Value = -4628.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(10, 8).Value = 7575
$ws.Cells.Item(10, 9).Value = 7920
$ws.Cells.Item(10, 10).Value = 7000
$ws.Cells.Item(10, 11).Value = 7920
$ws.Cells.Item(10, 12).Value = 7000
$ws.Cells.Item(10, 14).Value = -7340
$ws.Cells.Item(10, 13).Value = -7750
$ws.Cells.Item(55, 8).Value = 23990.166
$ws.Cells.Item(55, 9).Value = 0
$ws.Cells.Item(55, 10).Value = 23990.166
$ws.Cells.Item(55, 11).Value = 0
$ws.Cells.Item(55, 12).Value = 23990.166
$ws.Cells.Item(55, 13).ClearContents()
$ws.Cells.Item(55, 14).Value = -24620.166
$ws.Cells.Item(61, 8).Value = 2908.48
$ws.Cells.Item(61, 9).Value = 3048.2727
$ws.Cells.Item(61, 10).Value = 1883.3334
$ws.Cells.Item(61, 11).Value = 3048.2727
$ws.Cells.Item(61, 12).Value = 1883.3334
$ws.Cells.Item(61, 13).Value = -2836.2727
$ws.Cells.Item(61, 14).Value = -2307.3334
$ws.Cells.Item(74, 8).Value = 4912.2856
$ws.Cells.Item(74, 9).Value = 4972.75
$ws.Cells.Item(74, 11).Value = 4972.75
$ws.Cells.Item(74, 13).Value = -4098.75
$ws.Cells.Item(77, 8).Value = 4912.2856
$ws.Cells.Item(77, 9).Value = 4972.75
$ws.Cells.Item(77, 11).Value = 24863.75
$ws.Cells.Item(77, 13).Value = -20495.75
$ws.Cells.Item(122, 8).Value = 486332.1
$ws.Cells.Item(122, 9).Value = 672543.75
$ws.Cells.Item(122, 10).Value = 20803
$ws.Cells.Item(122, 11).Value = 2017631.25
$ws.Cells.Item(122, 12).Value = 62409
$ws.Cells.Item(122, 13).Value = -2015181.25
$ws.Cells.Item(122, 14).Value = -67309
$ws.Cells.Item(131, 8).Value = 80715
$ws.Cells.Item(131, 10).Value = 80715
$ws.Cells.Item(131, 12).Value = 80715
$ws.Cells.Item(131, 14).Value = -90795
$ws.Cells.Item(132, 8).Value = 4180.5
$ws.Cells.Item(132, 9).Value = 4861
$ws.Cells.Item(132, 10).Value = 3500
$ws.Cells.Item(132, 11).Value = 14583
$ws.Cells.Item(132, 12).Value = 10500
$ws.Cells.Item(132, 13).Value = -12053
$ws.Cells.Item(132, 14).Value = -15560
$ws.Cells.Item(136, 8).Value = 2908.48
$ws.Cells.Item(136, 9).Value = 3048.2727
$ws.Cells.Item(136, 10).Value = 1883.3334
$ws.Cells.Item(136, 11).Value = 9144.8181
$ws.Cells.Item(136, 12).Value = 5650.0002
$ws.Cells.Item(136, 13).Value = -6594.8181
$ws.Cells.Item(136, 14).Value = -10750.0002

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22, 8).Value = 194
$ws.Cells.Item(22, 9).Value = 196.85715
$ws.Cells.Item(22, 10).Value = 174
$ws.Cells.Item(22, 11).Value = 196.85715
$ws.Cells.Item(22, 12).Value = 174
$ws.Cells.Item(22, 13).Value = -23.85714999999999
$ws.Cells.Item(22, 14).Value = -520
$ws.Cells.Item(130, 8).Value = 59389.8
$ws.Cells.Item(130, 10).Value = 59389.8
$ws.Cells.Item(130, 12).Value = 59389.8
$ws.Cells.Item(130, 14).Value = -69429.8
$ws.Cells.Item(134, 8).Value = 3000
$ws.Cells.Item(134, 9).Value = 3000
$ws.Cells.Item(134, 11).Value = 9000
$ws.Cells.Item(134, 13).Value = -6465

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 4685.5
$ws.Cells.Item(58, 10).Value = 4685.5
$ws.Cells.Item(58, 12).Value = 4685.5
$ws.Cells.Item(58, 14).Value = -5091.5
$ws.Cells.Item(134, 8).Value = 3213.4285
$ws.Cells.Item(134, 9).Value = 1749.5
$ws.Cells.Item(134, 11).Value = 5248.5
$ws.Cells.Item(134, 13).Value = -2713.5
$ws.Cells.Item(136, 8).Value = 4685.5
$ws.Cells.Item(136, 10).Value = 4685.5
$ws.Cells.Item(136, 12).Value = 14056.5
$ws.Cells.Item(136, 14).Value = -19156.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(98, 8).Value = 2953.0908
$ws.Cells.Item(98, 10).Value = 2609.4443
$ws.Cells.Item(98, 12).Value = 7828.3329
$ws.Cells.Item(98, 14).Value = -10824.3329
$ws.Cells.Item(120, 8).Value = 14166.667
$ws.Cells.Item(120, 9).Value = 0
$ws.Cells.Item(120, 10).Value = 14166.667
$ws.Cells.Item(120, 11).Value = 0
$ws.Cells.Item(120, 12).Value = 42500.001
$ws.Cells.Item(120, 13).ClearContents()
$ws.Cells.Item(120, 14).Value = -52176.001

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(31, 8).Value = 1585.8
$ws.Cells.Item(31, 9).Value = 843.3333
$ws.Cells.Item(31, 10).Value = 2699.5
$ws.Cells.Item(31, 11).Value = 843.3333
$ws.Cells.Item(31, 12).Value = 2699.5
$ws.Cells.Item(31, 13).Value = -551.3333
$ws.Cells.Item(31, 14).Value = -3283.5
$ws.Cells.Item(37, 8).Value = 1585.8
$ws.Cells.Item(37, 9).Value = 843.3333
$ws.Cells.Item(37, 10).Value = 2699.5
$ws.Cells.Item(37, 11).Value = 843.3333
$ws.Cells.Item(37, 12).Value = 2699.5
$ws.Cells.Item(37, 13).Value = -566.3333
$ws.Cells.Item(37, 14).Value = -3253.5
$ws.Cells.Item(132, 8).Value = 1686.75
$ws.Cells.Item(132, 9).Value = 1686.75
$ws.Cells.Item(132, 11).Value = 5060.25
$ws.Cells.Item(132, 13).Value = -2530.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(4, 8).Value = 4785.6665
$ws.Cells.Item(4, 9).Value = 5603
$ws.Cells.Item(4, 10).Value = 3968.3333
$ws.Cells.Item(4, 11).Value = 5603
$ws.Cells.Item(4, 12).Value = 3968.3333
$ws.Cells.Item(4, 13).Value = -5490
$ws.Cells.Item(4, 14).Value = -4194.3333
$ws.Cells.Item(28, 8).Value = 4785.6665
$ws.Cells.Item(28, 9).Value = 5603
$ws.Cells.Item(28, 10).Value = 3968.3333
$ws.Cells.Item(28, 11).Value = 5603
$ws.Cells.Item(28, 12).Value = 3968.3333
$ws.Cells.Item(28, 13).Value = -5371
$ws.Cells.Item(28, 14).Value = -4432.3333
$ws.Cells.Item(37, 8).Value = 4785.6665
$ws.Cells.Item(37, 9).Value = 5603
$ws.Cells.Item(37, 10).Value = 3968.3333
$ws.Cells.Item(37, 11).Value = 5603
$ws.Cells.Item(37, 12).Value = 3968.3333
$ws.Cells.Item(37, 13).Value = -5496
$ws.Cells.Item(37, 14).Value = -4182.3333
$ws.Cells.Item(46, 8).Value = 2337.2
$ws.Cells.Item(46, 9).Value = 1834.125
$ws.Cells.Item(46, 11).Value = 1834.125
$ws.Cells.Item(46, 13).Value = -1646.125
$ws.Cells.Item(122, 8).Value = 4999
$ws.Cells.Item(122, 10).Value = 4999
$ws.Cells.Item(122, 12).Value = 14997
$ws.Cells.Item(122, 14).Value = -19897
$ws.Cells.Item(136, 8).Value = 5393.222
$ws.Cells.Item(136, 9).Value = 5079.4
$ws.Cells.Item(136, 11).Value = 15238.2
$ws.Cells.Item(136, 13).Value = -12688.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(38, 8).Value = 13056
$ws.Cells.Item(38, 9).Value = 13056
$ws.Cells.Item(38, 11).Value = 13056
$ws.Cells.Item(38, 13).Value = -12583
$ws.Cells.Item(41, 8).Value = 18333.334
$ws.Cells.Item(41, 10).Value = 18333.334
$ws.Cells.Item(41, 12).Value = 18333.334
$ws.Cells.Item(41, 14).Value = -19113.334
$ws.Cells.Item(45, 8).Value = 22539
$ws.Cells.Item(45, 10).Value = 26313
$ws.Cells.Item(45, 12).Value = 26313
$ws.Cells.Item(45, 14).Value = -27295
$ws.Cells.Item(48, 8).Value = 42350
$ws.Cells.Item(48, 9).Value = 40000
$ws.Cells.Item(48, 10).Value = 43525
$ws.Cells.Item(48, 11).Value = 40000
$ws.Cells.Item(48, 12).Value = 43525
$ws.Cells.Item(48, 13).Value = -39431
$ws.Cells.Item(48, 14).Value = -44663
$ws.Cells.Item(49, 8).Value = 264710.72
$ws.Cells.Item(49, 9).Value = 319663
$ws.Cells.Item(49, 10).Value = 223496.5
$ws.Cells.Item(49, 11).Value = 319663
$ws.Cells.Item(49, 12).Value = 223496.5
$ws.Cells.Item(49, 13).Value = -319433
$ws.Cells.Item(49, 14).Value = -223956.5
$ws.Cells.Item(107, 8).Value = 580.25
$ws.Cells.Item(107, 9).Value = 520.2857
$ws.Cells.Item(107, 11).Value = 1560.8571
$ws.Cells.Item(107, 13).Value = 359.1428999999998
$ws.Cells.Item(117, 8).Value = 62500
$ws.Cells.Item(117, 10).Value = 62500
$ws.Cells.Item(117, 12).Value = 62500
$ws.Cells.Item(117, 14).Value = -71678
$ws.Cells.Item(126, 8).Value = 1927.7142
$ws.Cells.Item(126, 9).Value = 1499.2
$ws.Cells.Item(126, 10).Value = 2999
$ws.Cells.Item(126, 11).Value = 4497.6
$ws.Cells.Item(126, 12).Value = 8997
$ws.Cells.Item(126, 13).Value = -2027.6
$ws.Cells.Item(126, 14).Value = -13937
$ws.Cells.Item(136, 8).Value = 4081.6924
$ws.Cells.Item(136, 9).Value = 3843.7778
$ws.Cells.Item(136, 11).Value = 11531.3334
$ws.Cells.Item(136, 13).Value = -8981.3334
